$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a value to E4 (Fabric column for Item 2002 / Fan Kit)
$ws.Range("E4").Value = 20

# Update the BOM reference for item 3000 (Standing Desk) to include the new item
$ws.Range("K5").Value = "2001,1;2002,1"

# Update the active selection to E3
$ws.Range("E3").Select()
